$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lama2"
$ws.Cells.Item(2, 3).Value = "Dag1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.471482
$ws.Cells.Item(2, 8).Value = 4.414446
$ws.Cells.Item(2, 9).Value = 0.004946458467382327
$ws.Cells.Item(2, 10).Value = 0.004946458467382326
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 11.779764
$ws.Cells.Item(2, 14).Value = 35.339292
$ws.Cells.Item(2, 15).Value = 0.1028447940505417
$ws.Cells.Item(2, 16).Value = 0.1028447940505417
$ws.Cells.Item(2, 17).Value = 17.333710690248
$ws.Cells.Item(2, 18).Value = 156.003396212232
$ws.Cells.Item(2, 19).Value = 0.0005087175023574936
$ws.Cells.Item(2, 20).Value = 0.0005087175023574935

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lama2"
$ws.Cells.Item(3, 3).Value = "Dag1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.471482
$ws.Cells.Item(3, 8).Value = 4.414446
$ws.Cells.Item(3, 9).Value = 0.004946458467382327
$ws.Cells.Item(3, 10).Value = 0.004946458467382326
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 43.841352
$ws.Cells.Item(3, 14).Value = 131.524056
$ws.Cells.Item(3, 15).Value = 0.3827627461243965
$ws.Cells.Item(3, 16).Value = 0.3827627461243964
$ws.Cells.Item(3, 17).Value = 64.51176032366399
$ws.Cells.Item(3, 18).Value = 580.605842912976
$ws.Cells.Item(3, 19).Value = 0.001893320026565533
$ws.Cells.Item(3, 20).Value = 0.001893320026565532

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lama2"
$ws.Cells.Item(4, 3).Value = "Dag1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.471482
$ws.Cells.Item(4, 8).Value = 4.414446
$ws.Cells.Item(4, 9).Value = 0.004946458467382327
$ws.Cells.Item(4, 10).Value = 0.004946458467382326
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 58.91811866666666
$ws.Cells.Item(4, 14).Value = 176.754356
$ws.Cells.Item(4, 15).Value = 0.5143924598250619
$ws.Cells.Item(4, 16).Value = 0.5143924598250619
$ws.Cells.Item(4, 17).Value = 86.69695109186399
$ws.Cells.Item(4, 18).Value = 780.2725598267759
$ws.Cells.Item(4, 19).Value = 0.002544420938459301
$ws.Cells.Item(4, 20).Value = 0.002544420938459301

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Lama2"
$ws.Cells.Item(5, 3).Value = "Dag1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 264.7713316666666
$ws.Cells.Item(5, 8).Value = 794.313995
$ws.Cells.Item(5, 9).Value = 0.8900417371348598
$ws.Cells.Item(5, 10).Value = 0.8900417371348596
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 11.779764
$ws.Cells.Item(5, 14).Value = 35.339292
$ws.Cells.Item(5, 15).Value = 0.1028447940505417
$ws.Cells.Item(5, 16).Value = 0.1028447940505417
$ws.Cells.Item(5, 17).Value = 3118.94380099906
$ws.Cells.Item(5, 18).Value = 28070.49420899154
$ws.Cells.Item(5, 19).Value = 0.09153615915202104
$ws.Cells.Item(5, 20).Value = 0.091536159152021

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lama2"
$ws.Cells.Item(6, 3).Value = "Dag1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 264.7713316666666
$ws.Cells.Item(6, 8).Value = 794.313995
$ws.Cells.Item(6, 9).Value = 0.8900417371348598
$ws.Cells.Item(6, 10).Value = 0.8900417371348596
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 43.841352
$ws.Cells.Item(6, 14).Value = 131.524056
$ws.Cells.Item(6, 15).Value = 0.3827627461243965
$ws.Cells.Item(6, 16).Value = 0.3827627461243964
$ws.Cells.Item(6, 17).Value = 11607.93315110708
$ws.Cells.Item(6, 18).Value = 104471.3983599637
$ws.Cells.Item(6, 19).Value = 0.3406748194710671
$ws.Cells.Item(6, 20).Value = 0.340674819471067

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Lama2"
$ws.Cells.Item(7, 3).Value = "Dag1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 264.7713316666666
$ws.Cells.Item(7, 8).Value = 794.313995
$ws.Cells.Item(7, 9).Value = 0.8900417371348598
$ws.Cells.Item(7, 10).Value = 0.8900417371348596
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 58.91811866666666
$ws.Cells.Item(7, 14).Value = 176.754356
$ws.Cells.Item(7, 15).Value = 0.5143924598250619
$ws.Cells.Item(7, 16).Value = 0.5143924598250619
$ws.Cells.Item(7, 17).Value = 15599.82873866802
$ws.Cells.Item(7, 18).Value = 140398.4586480122
$ws.Cells.Item(7, 19).Value = 0.4578307585117716
$ws.Cells.Item(7, 20).Value = 0.4578307585117715

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Lama2"
$ws.Cells.Item(8, 3).Value = "Dag1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 31.239114
$ws.Cells.Item(8, 8).Value = 93.717342
$ws.Cells.Item(8, 9).Value = 0.105011804397758
$ws.Cells.Item(8, 10).Value = 0.105011804397758
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 11.779764
$ws.Cells.Item(8, 14).Value = 35.339292
$ws.Cells.Item(8, 15).Value = 0.1028447940505417
$ws.Cells.Item(8, 16).Value = 0.1028447940505417
$ws.Cells.Item(8, 17).Value = 367.989390489096
$ws.Cells.Item(8, 18).Value = 3311.904514401864
$ws.Cells.Item(8, 19).Value = 0.0107999173961632
$ws.Cells.Item(8, 20).Value = 0.01079991739616319

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Lama2"
$ws.Cells.Item(9, 3).Value = "Dag1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 31.239114
$ws.Cells.Item(9, 8).Value = 93.717342
$ws.Cells.Item(9, 9).Value = 0.105011804397758
$ws.Cells.Item(9, 10).Value = 0.105011804397758
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 43.841352
$ws.Cells.Item(9, 14).Value = 131.524056
$ws.Cells.Item(9, 15).Value = 0.3827627461243965
$ws.Cells.Item(9, 16).Value = 0.3827627461243964
$ws.Cells.Item(9, 17).Value = 1369.564993042128
$ws.Cells.Item(9, 18).Value = 12326.08493737915
$ws.Cells.Item(9, 19).Value = 0.04019460662676384
$ws.Cells.Item(9, 20).Value = 0.04019460662676383

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Lama2"
$ws.Cells.Item(10, 3).Value = "Dag1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 31.239114
$ws.Cells.Item(10, 8).Value = 93.717342
$ws.Cells.Item(10, 9).Value = 0.105011804397758
$ws.Cells.Item(10, 10).Value = 0.105011804397758
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 58.91811866666666
$ws.Cells.Item(10, 14).Value = 176.754356
$ws.Cells.Item(10, 15).Value = 0.5143924598250619
$ws.Cells.Item(10, 16).Value = 0.5143924598250619
$ws.Cells.Item(10, 17).Value = 1840.549825693528
$ws.Cells.Item(10, 18).Value = 16564.94843124175
$ws.Cells.Item(10, 19).Value = 0.05401728037483101
$ws.Cells.Item(10, 20).Value = 0.054017280374831
